# "Generate Report for Archive"
#
# Localization status moved on from "Ready for handoff" -> "In Translation"
# for the e2e/*.md and e2e/*.png rows, so every cell that showed the old
# status text is updated on the Overview sheet (columns zh-cn/de-de) and on
# each per-locale detail sheet (Status column).
#
# Excel re-measured ("best fit") the Status-ish columns after the text
# changed (the new text is shorter than the old), which is why the report's
# column widths shrink along with the text edit.

$wb  = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status columns, rows 2-4 ---
foreach ($row in 2..4) {
    foreach ($col in @("E", "F")) {
        $cell = $overview.Range("$col$row")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Per-locale detail sheets: Status column (C), rows 2-4 ---
foreach ($ws in @($zhcn, $dede)) {
    foreach ($row in 2..4) {
        $cell = $ws.Range("C$row")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Column widths: re-fit the columns that held the status text ---
# (Excel stores column width as a character count on a 1/MDW pixel grid, so
# the setter below snaps to the nearest representable width for the target
# pixel run - same behaviour as typing a width into the Format > Column
# Width dialog.)
$overview.Range("E1").ColumnWidth = 16.334635416666664
$overview.Range("F1").ColumnWidth = 16.334635416666664
$zhcn.Range("C1").ColumnWidth     = 12.501302083333332
$dede.Range("C1").ColumnWidth     = 12.501302083333332
